# Fix bug in significance analysis: the "num_XXX" columns (E,G,I,K,M,O)
# were incorrectly defaulting to 5 in many rows. This script restores the
# correct significant-cell counts and recomputes the dependent p_XXX
# columns (D,F,H,J,L,N) as num_XXX / p_cells (column C), matching the
# original data pipeline.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Corrected "num_XXX" (significant cell count) values, keyed by cell address.
$numFixes = @{
    "E3" = 2
    "G3" = 2
    "M3" = 1
    "O3" = 2
    "I5" = 1
    "K5" = 1
    "M5" = 2
    "E6" = 1
    "G6" = 3
    "K6" = 1
    "M6" = 3
    "O6" = 4
    "E7" = 1
    "G7" = 1
    "I7" = 2
    "K7" = 1
    "M7" = 1
    "O7" = 1
    "E8" = 3
    "G8" = 3
    "I8" = 4
    "M8" = 4
    "E9" = 1
    "G9" = 3
    "I9" = 4
    "M9" = 2
    "E10" = 1
    "G10" = 1
    "I10" = 2
    "M10" = 2
    "E11" = 2
    "G11" = 3
    "I11" = 4
    "O11" = 1
    "E12" = 3
    "G12" = 3
    "K12" = 3
    "M12" = 3
    "E13" = 3
    "M13" = 1
    "O13" = 4
    "E14" = 2
    "G14" = 3
    "I14" = 3
    "M14" = 1
    "O14" = 2
    "E15" = 1
    "G15" = 2
    "I15" = 3
    "O15" = 1
    "E16" = 2
    "G16" = 2
    "I16" = 4
    "K16" = 1
    "M16" = 1
    "O16" = 1
    "E19" = 3
    "G19" = 4
    "O19" = 3
    "E20" = 1
    "G20" = 2
    "I20" = 2
    "K20" = 1
    "M20" = 3
    "O20" = 3
    "I22" = 1
    "M22" = 2
    "E24" = 1
    "G24" = 3
    "I24" = 3
    "K24" = 1
    "M24" = 1
    "O24" = 1
    "K29" = 2
    "M29" = 4
    "E32" = 1
    "O32" = 2
    "E37" = 1
    "G37" = 1
    "I37" = 1
    "K37" = 3
    "K38" = 1
    "M38" = 1
    "O38" = 2
    "E41" = 1
    "G41" = 2
    "I41" = 3
    "M41" = 1
    "O41" = 2
    "E44" = 4
    "G44" = 4
    "K44" = 1
    "M44" = 4
    "E45" = 3
    "G45" = 3
    "I45" = 4
    "E46" = 1
    "G46" = 3
    "K46" = 4
    "E47" = 1
    "G47" = 1
    "I47" = 2
    "M47" = 1
    "O47" = 2
    "E48" = 3
    "G48" = 4
    "K48" = 2
    "M48" = 2
    "O48" = 4
    "G49" = 2
    "I49" = 2
    "O49" = 4
    "G50" = 1
    "I50" = 1
    "K50" = 3
    "M50" = 4
    "K51" = 1
    "E52" = 1
    "G52" = 1
    "I52" = 3
    "K52" = 1
    "M52" = 1
    "O52" = 3
    "E54" = 2
    "G54" = 2
    "I54" = 2
    "M54" = 1
    "O54" = 4
}

foreach ($addr in $numFixes.Keys) {
    $ws.Range($addr).Value2 = $numFixes[$addr]
}

# Column letters: num_025=E, num_05=G, num_10=I, num_975=K, num_95=M, num_90=O
# paired with: p_025=D,  p_05=F,  p_10=H,  p_975=J, p_95=L,  p_90=N
$pairs = @{
    "D" = "E"
    "F" = "G"
    "H" = "I"
    "J" = "K"
    "L" = "M"
    "N" = "O"
}

for ($row = 2; $row -le 55; $row++) {
    $pCells = $ws.Cells.Item($row, 3).Value2
    if ($pCells -eq $null -or $pCells -eq 0) {
        continue
    }
    foreach ($pCol in $pairs.Keys) {
        $numCol = $pairs[$pCol]
        $numAddr = "$numCol$row"
        $pAddr = "$pCol$row"
        $numVal = $ws.Range($numAddr).Value2
        if ($numVal -eq $null) {
            continue
        }
        $ws.Range($pAddr).Value2 = $numVal / $pCells
    }
}
